# Rerun/summarise models without urban landuse:
#  - rename each summary sheet to its new run id
#  - update the "Education[T.Unknown]" label to "Education[T.Unknown/Other]"
#    on row 5 of every sheet

$wb = $excel.ActiveWorkbook

$oldNames = @(
    "summ40031053",
    "summ40329361",
    "summ40621792",
    "summ40871029",
    "summ41136319",
    "summ41438727",
    "summ41717010",
    "summ42076757",
    "summ42367286"
)

$newNames = @(
    "summ24812729",
    "summ25164278",
    "summ25443093",
    "summ25712460",
    "summ25988158",
    "summ26239703",
    "summ26498173",
    "summ26768534",
    "summ27041379"
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)

    if ($ws.Name -eq $oldNames[$i]) {
        $ws.Name = $newNames[$i]
    }

    if ($ws.Range("A5").Value2 -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
